$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1088.36
$ws.Range("I15").Value = 1088.36
$ws.Range("K15").Value = 3265.08
$ws.Range("M15").Value = -3096.08
$ws.Range("H40").Value = 875
$ws.Range("I40").Value = 875
$ws.Range("K40").Value = 875
$ws.Range("M40").Value = -700
$ws.Range("H137").Value = 2274.22
$ws.Range("I137").Value = 1646.48
$ws.Range("J137").Value = 2901.96
$ws.Range("K137").Value = 4939.440000000001
$ws.Range("L137").Value = 8705.880000000001
$ws.Range("M137").Value = -2389.440000000001
$ws.Range("N137").Value = -13805.88
$ws.Range("H138").Value = 2441.058
$ws.Range("I138").Value = 1630.9445
$ws.Range("J138").Value = 2655.5
$ws.Range("K138").Value = 4892.833500000001
$ws.Range("L138").Value = 7966.5
$ws.Range("M138").Value = 247.1664999999994
$ws.Range("N138").Value = -18246.5
$ws.Range("H140").Value = 37495
$ws.Range("J140").Value = 37495
$ws.Range("L140").Value = 37495
$ws.Range("N140").Value = -47855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("N9").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
$ws.Range("H32").Value = 10212
$ws.Range("I32").Value = 7532.5615
$ws.Range("K32").Value = 7532.5615
$ws.Range("M32").Value = -7245.5615
$ws.Range("H37").Value = 18444
$ws.Range("I37").Value = 8888
$ws.Range("J37").Value = 28000
$ws.Range("K37").Value = 8888
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = -8615
$ws.Range("N37").Value = -28546
$ws.Range("H44").Value = 29800
$ws.Range("J44").Value = 29800
$ws.Range("L44").Value = 29800
$ws.Range("N44").Value = -30776
$ws.Range("H45").Value = 1119
$ws.Range("I45").Value = 872.0909
$ws.Range("J45").Value = 1571.6666
$ws.Range("K45").Value = 872.0909
$ws.Range("L45").Value = 1571.6666
$ws.Range("M45").Value = -495.0909
$ws.Range("N45").Value = -2325.6666
$ws.Range("H55").Value = 38000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H74").Value = 1366.1714
$ws.Range("I74").Value = 860.8148
$ws.Range("J74").Value = 3071.75
$ws.Range("K74").Value = 860.8148
$ws.Range("L74").Value = 3071.75
$ws.Range("M74").Value = 13.18520000000001
$ws.Range("N74").Value = -4819.75
$ws.Range("H77").Value = 1366.1714
$ws.Range("I77").Value = 860.8148
$ws.Range("J77").Value = 3071.75
$ws.Range("K77").Value = 4304.074
$ws.Range("L77").Value = 15358.75
$ws.Range("M77").Value = 63.92600000000039
$ws.Range("N77").Value = -24094.75
$ws.Range("H132").Value = 4571.143
$ws.Range("I132").Value = 4721.6924
$ws.Range("K132").Value = 14165.0772
$ws.Range("M132").Value = -11635.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 1452.5
$ws.Range("I10").Value = 1452.5
$ws.Range("K10").Value = 1452.5
$ws.Range("M10").Value = -1312.5
$ws.Range("H22").Value = 484.625
$ws.Range("I22").Value = 466.66666
$ws.Range("J22").Value = 495.4
$ws.Range("K22").Value = 466.66666
$ws.Range("L22").Value = 495.4
$ws.Range("M22").Value = -293.66666
$ws.Range("N22").Value = -841.4
$ws.Range("H24").Value = 408
$ws.Range("I24").Value = 316
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 316
$ws.Range("L24").Value = 500
$ws.Range("M24").Value = -81
$ws.Range("N24").Value = -970
$ws.Range("H99").Value = 100001230
$ws.Range("I99").Value = 142858190
$ws.Range("K99").Value = 142858190
$ws.Range("M99").Value = -142856692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1607.2307
$ws.Range("I31").Value = 1403.8096
$ws.Range("J31").Value = 2461.6
$ws.Range("K31").Value = 1403.8096
$ws.Range("L31").Value = 2461.6
$ws.Range("M31").Value = -1108.8096
$ws.Range("N31").Value = -3051.6
$ws.Range("H34").Value = 1607.2307
$ws.Range("I34").Value = 1403.8096
$ws.Range("J34").Value = 2461.6
$ws.Range("K34").Value = 1403.8096
$ws.Range("L34").Value = 2461.6
$ws.Range("M34").Value = -1201.8096
$ws.Range("N34").Value = -2865.6
$ws.Range("H58").Value = 6041.7393
$ws.Range("I58").Value = 1204.5
$ws.Range("J58").Value = 13566.333
$ws.Range("K58").Value = 1204.5
$ws.Range("L58").Value = 13566.333
$ws.Range("M58").Value = -1001.5
$ws.Range("N58").Value = -13972.333
$ws.Range("H62").Value = 15386793
$ws.Range("I62").Value = 2386.4285
$ws.Range("J62").Value = 33335266
$ws.Range("K62").Value = 2386.4285
$ws.Range("L62").Value = 33335266
$ws.Range("M62").Value = -1762.4285
$ws.Range("N62").Value = -33336514
$ws.Range("H65").Value = 15386793
$ws.Range("I65").Value = 2386.4285
$ws.Range("J65").Value = 33335266
$ws.Range("K65").Value = 11932.1425
$ws.Range("L65").Value = 166676330
$ws.Range("M65").Value = -8812.1425
$ws.Range("N65").Value = -166682570
$ws.Range("H136").Value = 6041.7393
$ws.Range("I136").Value = 1204.5
$ws.Range("J136").Value = 13566.333
$ws.Range("K136").Value = 3613.5
$ws.Range("L136").Value = 40698.999
$ws.Range("M136").Value = -1063.5
$ws.Range("N136").Value = -45798.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 13994.389
$ws.Range("I3").Value = 7981
$ws.Range("J3").Value = 20007.777
$ws.Range("K3").Value = 23943
$ws.Range("L3").Value = 60023.33099999999
$ws.Range("M3").Value = -23831
$ws.Range("N3").Value = -60247.33099999999
$ws.Range("H34").Value = 4547225.5
$ws.Range("J34").Value = 7145441.5
$ws.Range("L34").Value = 21436324.5
$ws.Range("N34").Value = -21436492.5
$ws.Range("H39").Value = 4284.923
$ws.Range("J39").Value = 4373.091
$ws.Range("L39").Value = 13119.273
$ws.Range("N39").Value = -13707.273
$ws.Range("H55").Value = 3250
$ws.Range("J55").Value = 3250
$ws.Range("L55").Value = 9750
$ws.Range("N55").Value = -10104
$ws.Range("H74").Value = 4500
$ws.Range("J74").Value = 4500
$ws.Range("L74").Value = 13500
$ws.Range("N74").Value = -15622
$ws.Range("H77").Value = 4500
$ws.Range("J77").Value = 4500
$ws.Range("L77").Value = 40500
$ws.Range("N77").Value = -51108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 4900
$ws.Range("J19").Value = 4900
$ws.Range("L19").Value = 4900
$ws.Range("N19").Value = -5476
$ws.Range("H95").Value = 19299.666
$ws.Range("J95").Value = 19299.666
$ws.Range("L95").Value = 19299.666
$ws.Range("N95").Value = -24791.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 226.33333
$ws.Range("I55").Value = 179.84616
$ws.Range("K55").Value = 179.84616
$ws.Range("M55").Value = -6.846159999999998
$ws.Range("H100").Value = 1976
$ws.Range("I100").Value = 1990
$ws.Range("J100").Value = 1966.6666
$ws.Range("K100").Value = 1990
$ws.Range("L100").Value = 1966.6666
$ws.Range("M100").Value = -1449
$ws.Range("N100").Value = -3048.6666
$ws.Range("H132").Value = 2571.5715
$ws.Range("I132").Value = 2146.6155
$ws.Range("J132").Value = 2939.8667
$ws.Range("K132").Value = 6439.8465
$ws.Range("L132").Value = 8819.6001
$ws.Range("M132").Value = -3909.8465
$ws.Range("N132").Value = -13879.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40713.57
$ws.Range("J46").Value = 40713.57
$ws.Range("L46").Value = 40713.57
$ws.Range("N46").Value = -41175.57
$ws.Range("H122").Value = 8623002
$ws.Range("I122").Value = 12502386
$ws.Range("K122").Value = 37507158
$ws.Range("M122").Value = -37504708
$ws.Range("H134").Value = 40713.57
$ws.Range("J134").Value = 40713.57
$ws.Range("L134").Value = 122140.71
$ws.Range("N134").Value = -127210.71
